# "Added last minute updates"
#
# The first paragraph of the body currently reads:
#     **ID__AFFARS_5350_topic_8__ID**<space>
# (two runs: the placeholder text, then a run holding a single trailing
# space) with pPr: <w:spacing w:after="0"/><w:ind w:left="120"/><w:jc w:val="left"/>
#
# It needs to become a single run reading:
#     **ID__AFFARS_5350_103_5__ID**
# (no trailing space / no second run) and the paragraph's pPr should gain a
# paragraph border (5-twip gap on all four sides, matching the borders
# already used lower in the document) and its left indent should grow from
# 120 twips (6pt) to 225 twips (11.25pt) -- again matching the indent used
# by the already-bordered paragraphs further down in the document.

$d = $word.ActiveDocument

# Replace the old topic id placeholder (plus the trailing-space run that
# immediately follows it) with the new placeholder text. Because both runs
# share identical run formatting, Word's Find/Replace collapses them into a
# single run and drops the now-absent trailing whitespace.
$d.Content.Find.Execute(
    "**ID__AFFARS_5350_topic_8__ID** ", $true, $false, $false, $false, $false,
    $true, 1, $false, "**ID__AFFARS_5350_103_5__ID**", 2)

# The updated text lives in the first paragraph of the document.
$p1 = $d.Paragraphs(1)

# Grow the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Add the (style-less) paragraph border -- 5-twip spacing on every edge,
# same as the pattern already used by the paragraphs below.
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
